$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "MCT-2A-EAP"
$ws.Range("F3").Value = "MEC-3A-EAP"

$ws.Range("E4").Value = "MCT-2A-EAP"
$ws.Range("F4").Value = "MEC-3A-EAP"
